# Remove log rows that are not for "data11.xlsx" (delete functionality)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows from bottom to top so row indices of the rows
# still to be removed are unaffected by the shift caused by deletion.
$ws.Rows.Item(12).Delete()   # 65bddbed... | data14.xlsx
$ws.Rows.Item(11).Delete()   # d465c3ad... | data13.xlsx
$ws.Rows.Item(10).Delete()   # 8fdc9d60... | data11 - Copy.xlsx
$ws.Rows.Item(7).Delete()    # 89597334... | data11 - Copy.xlsx
$ws.Rows.Item(3).Delete()    # eae00c4f... | data11 - Copy.xlsx
